# Netflix.xlsx edit: "Rimosso il 2022 + .gitignore"
#
# Semantics of the change (reverse-engineered from the target diff):
#   The sheet tracks one row per month across many years (row 1 = month
#   number 1-12, row 5 = year, rows 2/3/4 = ratios derived from rows 8/9).
#   The whole "year 2022" (the first 12 month-columns, B:M) is removed by
#   shifting the header/ratio rows (2, 3, 4, 5 -- NOT row 1, which is a
#   pure self-referential 12-month cycle and doesn't need to move) left by
#   12 columns: column X now shows what used to be at column X+12. The
#   tail 12 columns that no longer have source data (BV:CG) become empty.
#   Row 5's own running formula keeps referencing "row 1, 12 columns to
#   the right" (un-shifted, since row 1 never moved) while its
#   self-reference (previous cell in row 5) follows the shift normally.
#
# Additionally the active sheet view resets its frozen/scrolled
# topLeftCell back to A1 and the selection moves to G13, and an empty
# (but centered-style) cell is introduced at G13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColLetter([int]$n) {
    $letters = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letters = [char](65 + $rem) + $letters
        $n = [int](($n - $rem - 1) / 26)
    }
    return $letters
}

$firstDataCol = 2    # B   -- first month column
$lastDataCol  = 85   # CG  -- last month column
$shift        = 12   # one year of months being removed

# ---------------------------------------------------------------------
# Row 2: ratio cells only change from column O (15) onward; B2:N2 keep
# their existing "checkmark" shared-string cells untouched.
#   new[c] = old[c+shift]   for c in O(15)..BU(73)
#   column O itself is a hard-coded number (0.8); P..BU are formulas
#   that reference row 8/row 9 at their OWN (old) column, so the
#   formula text is simply copied across unchanged.
# ---------------------------------------------------------------------
$startC = 15   # O
$endC   = 73   # BU
$n = $endC - $startC + 1
$arr = New-Object 'object[,]' 1,$n
for ($i = 0; $i -lt $n; $i++) {
    $c = $startC + $i
    $srcCol = ColLetter($c + $shift)
    if ($c -eq $startC) {
        $arr[0,$i] = 0.8
    } else {
        $arr[0,$i] = "=ROUND(" + $srcCol + "8/" + $srcCol + "9,  1)"
    }
}
$rng = (ColLetter $startC) + "2:" + (ColLetter $endC) + "2"
$ws.Range($rng).Formula = $arr

# ---------------------------------------------------------------------
# Row 3 and Row 4: identical pattern. Every column from B(2) to BU(73)
# becomes =ROUND(X8/X9, 1) where X is the OLD column (c+shift) -- i.e.
# the same ratio formula that used to live 12 columns to the right.
# ---------------------------------------------------------------------
$startC = 2    # B
$endC   = 73   # BU
$n = $endC - $startC + 1
$arr3 = New-Object 'object[,]' 1,$n
for ($i = 0; $i -lt $n; $i++) {
    $c = $startC + $i
    $srcCol = ColLetter($c + $shift)
    $arr3[0,$i] = "=ROUND(" + $srcCol + "8/" + $srcCol + "9,  1)"
}
$rng3 = (ColLetter $startC) + "3:" + (ColLetter $endC) + "3"
$ws.Range($rng3).Formula = $arr3
$rng4 = (ColLetter $startC) + "4:" + (ColLetter $endC) + "4"
$ws.Range($rng4).Formula = $arr3

# ---------------------------------------------------------------------
# Row 5: running year number.
#   B5 = 2023 (hard-coded seed, replacing the old 2022)
#   C5..BU5 = IF(<col+11>1=12, <col-1>5+1, <col-1>5)
#     -- the row-1 reference stays 12 columns further right than the
#        row-5 self reference because row 1 itself never shifted while
#        row 5 did.
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 2023

$startC = 3    # C
$endC   = 73   # BU
$n = $endC - $startC + 1
$arr5 = New-Object 'object[,]' 1,$n
for ($i = 0; $i -lt $n; $i++) {
    $c = $startC + $i
    $row1Ref = ColLetter($c + $shift - 1)
    $row5Ref = ColLetter($c - 1)
    $arr5[0,$i] = "=IF(" + $row1Ref + "1=12," + $row5Ref + "5+1," + $row5Ref + "5)"
}
$rng5 = (ColLetter $startC) + "5:" + (ColLetter $endC) + "5"
$ws.Range($rng5).Formula = $arr5

# ---------------------------------------------------------------------
# Clear the trailing 12 columns (BV:CG) in rows 2-5: there is no more
# source data to shift into them now that a whole year was removed.
# ---------------------------------------------------------------------
$tailStart = 74  # BV
$tailEnd   = 85  # CG
$tailRange = (ColLetter $tailStart) + "2:" + (ColLetter $tailEnd) + "5"
$ws.Range($tailRange).ClearContents()

# ---------------------------------------------------------------------
# New (empty) styled cell at G13, matching the centered style used by
# other label cells such as A8/A9.
# ---------------------------------------------------------------------
$g13 = $ws.Range("G13")
$g13.HorizontalAlignment = -4108   # xlCenter
$g13.VerticalAlignment   = -4108   # xlCenter

# ---------------------------------------------------------------------
# View: drop the frozen/scrolled topLeftCell (back to A1) and move the
# selection to G13.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$g13.Select()

$excel.Calculate()
